$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '68.066.70'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').Value = '3.674.55'
$ws.Range('E3').Value = '  -3.12%  '
$ws.Range('E4').Value = '  -0.09%  '
Set-TextValue 'D5' '596.26'
$ws.Range('E5').Value = '  +0.43%  '
Set-TextValue 'D6' '165.79'
$ws.Range('E6').Value = '  -3.80%  '
$ws.Range('D7').Value = '3.672.60'
$ws.Range('E7').Value = '  -3.18%  '
$ws.Range('E8').Value = '  +0.00%  '
Set-TextValue 'D9' '0.532'
$ws.Range('E9').Value = '  +0.82%  '
$ws.Range('E10').Value = '  +3.32%  '
Set-TextValue 'D11' '6.26'
$ws.Range('E11').Value = '  -0.21%  '
Set-TextValue 'D12' '0.456'
$ws.Range('E12').Value = '  -1.71%  '
$ws.Range('E13').Value = '  -0.65%  '
$ws.Range('E14').Value = '  -0.20%  '
$ws.Range('D15').Value = '4.286.89'
$ws.Range('E15').Value = '  -3.19%  '
$ws.Range('D16').Value = '3.677.64'
$ws.Range('E16').Value = '  -3.09%  '
$ws.Range('D17').Value = '68.089.71'
$ws.Range('E17').Value = '  +0.05%  '
Set-TextValue 'D18' '7.21'
$ws.Range('E18').Value = '  +0.85%  '
Set-TextValue 'D19' '0.115'
$ws.Range('E19').Value = '  -0.95%  '
$ws.Range('E20').Value = '  +6.53%  '
Set-TextValue 'D21' '489.27'
$ws.Range('E21').Value = '  +0.15%  '
Set-TextValue 'D22' '9.05'
$ws.Range('E22').Value = '  -2.10%  '
Set-TextValue 'D23' '0.719'
$ws.Range('E23').Value = '  -1.73%  '
Set-TextValue 'D24' '84.28'
$ws.Range('E24').Value = '  -0.49%  '
$ws.Range('E25').Value = '  +3.52%  '
Set-TextValue 'D26' '2.27'
$ws.Range('E26').Value = '  -4.12%  '
Set-TextValue 'D27' '12.14'
$ws.Range('E27').Value = '  -0.47%  '
Set-TextValue 'D28' '10.01'
$ws.Range('E28').Value = '  -1.72%  '
Set-TextValue 'D29' '1.00'
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('E30').Value = '  -0.86%  '
Set-TextValue 'D31' '2.38'
$ws.Range('E31').Value = '  -2.12%  '
$ws.Range('E32').Value = '  +1.46%  '
Set-TextValue 'D33' '31.20'
$ws.Range('E33').Value = '  -4.66%  '
$ws.Range('D34').Value = '3.813.63'
$ws.Range('E34').Value = '  -3.10%  '
$ws.Range('E35').Value = '  -1.61%  '
$ws.Range('D36').Value = '3.615.57'
$ws.Range('E36').Value = '  -3.12%  '
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('E38').Value = '  -1.40%  '
$ws.Range('E39').Value = '  -0.67%  '
$ws.Range('E40').Value = '  -3.66%  '
$ws.Range('E41').Value = '  -1.61%  '
Set-TextValue 'D42' '432.45'
$ws.Range('E42').Value = '  -4.73%  '
Set-TextValue 'D43' '48.69'
$ws.Range('E43').Value = '  -0.44%  '
Set-TextValue 'D44' '1.94'
$ws.Range('E44').Value = '  -2.43%  '
$ws.Range('E45').Value = '  -2.70%  '
$ws.Range('E46').Value = '  +1.11%  '
$ws.Range('E47').Value = '  +0.00%  '
Set-TextValue 'D48' '40.27'
$ws.Range('E48').Value = '  -2.85%  '
Set-TextValue 'D49' '141.36'
$ws.Range('E49').Value = '  +1.52%  '
$ws.Range('D50').Value = '2.726.98'
$ws.Range('E50').Value = '  -3.71%  '
Set-TextValue 'D51' '0.0347'
$ws.Range('E51').Value = '  -1.12%  '
